# 02项目计划表.xlsx -- apply "Add files via upload" edit
#
# Summary of the change:
#  - Fill in completion percentages (column C) for the 2018.11.14 week block
#    (rows 173-178) and update a couple of task descriptions.
#  - Duplicate that week's block to create a new week block
#    (日期：2018.11.15 第十一周周四) in rows 181-190, with updated task text.
#  - Minor cosmetic view changes (selected cell, tab ratio).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Duplicate the previous week's block (rows 171:180) down into the new
#    block (rows 181:190) BEFORE editing any of the source text, so the
#    duplicate keeps the original (unedited) wording -- exactly like a
#    user would do by selecting the block, copying it, and pasting it
#    below before typing the new week's updates.
# ---------------------------------------------------------------------
$ws.Range("A171:D180").Copy()
$ws.Range("A181:D190").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A171:D180").Copy()
$ws.Range("A181:D190").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Update the existing (2018.11.14) week block, rows 173-178:
#    fill in the "完成情况" (completion %) column, and refresh a couple of
#    task descriptions.
# ---------------------------------------------------------------------

# 黄成志's task description changed.
$ws.Range("B174").Value = "初步实现后台注册代码"

# 郑嘉蔚 was on leave this week -- clear her task/percent, note it in 备注.
$ws.Range("B176").Value = ""
$ws.Range("D176").Value = "请假"

# Completion percentages.
$ws.Range("C173").Value = 0.5
$ws.Range("C174").Value = 0.8
$ws.Range("C175").Value = 0.8
$ws.Range("C177").Value = 0.5
$ws.Range("C178").Value = 0.5

$ws.Range("C173,C174,C175,C176,C177,C178").NumberFormat = "0%"

# Week summary text.
$ws.Range("A179").Value = "总结：消息接口完成，继续接入通知信息接口。"

# ---------------------------------------------------------------------
# 3. Fill in the new (2018.11.15) week block, rows 181-190.
# ---------------------------------------------------------------------

# 练富珊's updated task.
$ws.Range("B183").Value = "接入环信通知接口"

# Date header for the new week (set last, matching creation order of the
# new shared strings).
$ws.Range("A181").Value = "日期：2018.11.15 第十一周周四"

# Give the new week's percent-complete column the same percent format as
# the week above (values are left blank, matching the template row).
$ws.Range("C183,C184,C185,C186,C187,C188").NumberFormat = "0%"

# ---------------------------------------------------------------------
# 4. Merged cells for the two new "banner" rows.
# ---------------------------------------------------------------------
$ws.Range("A181:D181").Merge()
$ws.Range("A189:D190").Merge()

# ---------------------------------------------------------------------
# 5. Cosmetic view state (best effort).
# ---------------------------------------------------------------------
$ws.Range("H176").Select()
$excel.ActiveWindow.ScrollRow = 163
$excel.ActiveWindow.TabRatio = 0.585714285714286
